$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Center of Gravity to Rear Axle" (row 8) and
# "Center of Gravity to Front Axle" (row 9) variable rows.
# This shifts the Gs of Acceleration / Gs of Deceleration rows
# up from rows 10-11 to rows 8-9.
$ws.Rows("8:9").Delete()

# Update the active selection to match the new sheet extent.
$ws.Range("A12").Select()
